$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new consolidated row values. Each card's multiple attribute rows
# are now collapsed into a single Python-tuple-formatted string per card.
$colossalWhale = "('Colossal Whale', ['{5}{U}{U}', 'Creature " + [char]0x2014 + " Whale', 'Islandwalk (This creature can" + [char]0x2019 + "t be blocked as long as defending player controls an Island.)', 'Whenever Colossal Whale attacks, you may exile target creature defending player controls until Colossal Whale leaves the battlefield. (That creature returns under its owner" + [char]0x2019 + "s control.)', '5/5'])"

$goblinDiplomats = "('Goblin Diplomats', ['{1}{R}', 'Creature " + [char]0x2014 + " Goblin', '{T}: Each creature attacks this turn if able.', '2/1'])"

$hiveStirrings = "('Hive Stirrings', ['{2}{W}', 'Sorcery', 'Create two 1/1 colorless Sliver creature tokens.'])"

$meganticSliver = "('Megantic Sliver', ['{5}{G}', 'Creature " + [char]0x2014 + " Sliver', 'Sliver creatures you control get +3/+3.', '3/3'])"

$ratchetBomb = "('Ratchet Bomb', ['{2}', 'Artifact', '{T}: Put a charge counter on Ratchet Bomb.', '{T}, Sacrifice Ratchet Bomb: Destroy each nonland permanent with converted mana cost equal to the number of charge counters on Ratchet Bomb.'])"

# Clear out the old rows 2-26 first, then write the new consolidated rows.
$ws.Range("A2:A26").ClearContents()

$ws.Range("A2").Value = $colossalWhale
$ws.Range("A3").Value = $goblinDiplomats
$ws.Range("A4").Value = $hiveStirrings
$ws.Range("A5").Value = $meganticSliver
$ws.Range("A6").Value = $ratchetBomb
